# Update the "Summary" worksheet:
#  - Add a new "Three reactions charring" block in columns E:F (mirrors the
#    existing "Two reactions charring" block in C:D)
#  - Add new rows under "Two reactions charring" (D5 entry + a whole new
#    "1 Mass + TGA" / "2 Mass + TGA" row pair)
#  - Fix the "Single reaction non-charring with ignition" table: B15 was "?"
#    and should be "O", plus two new rows (1 Mass + TGA / 1 Mass + F Temp +
#    B Temp) get added under it.
#  - The "X/O" label used in D3/D4 is replaced by the clearer
#    "10 kW X/100 kW O" (site-wide rename of that shared string).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# ---- Two reactions charring (C:D) updates -------------------------------
# D3/D4 previously held "X/O" -> now "10 kW X/100 kW O"
$ws.Range("D3").Value2 = "10 kW X/100 kW O"
$ws.Range("D4").Value2 = "10 kW X/100 kW O"

# New row 5 entry for this block
$ws.Range("D5").Value2 = "X"
$ws.Range("D5").HorizontalAlignment = -4108  # xlCenter

# New row 6: "1 Mass + TGA" / "10 kW X/100 kW O"
$ws.Range("C6").Value2 = "1 Mass + TGA"
$ws.Range("D6").Value2 = "10 kW X/100 kW O"
$ws.Range("D6").HorizontalAlignment = -4108  # xlCenter

# New row 7: "2 Mass + TGA"
$ws.Range("C7").Value2 = "2 Mass + TGA"

# ---- New "Three reactions charring" block (E:F) --------------------------
$ws.Range("E1").Value2 = "Three reactions charring"
$ws.Range("E1:F1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("E1:F1").Merge()

$ws.Range("E2").Value2 = "Scheme"
$ws.Range("F2").Value2 = "Successful?"

$ws.Range("E3").Value2 = "1 Mass + TGA + DSC"
$ws.Range("E4").Value2 = "2 Mass"
$ws.Range("E5").Value2 = "1 Mass + F Temp + B Temp"
$ws.Range("E6").Value2 = "1 Mass + TGA"

# F3:F6 stay blank but carry the same centered style used across column D
$ws.Range("F3:F6").HorizontalAlignment = -4108  # xlCenter

# ---- "Single reaction non-charring with ignition" table fix -------------
$ws.Range("B15").Value2 = "O"
$ws.Range("A16").Value2 = "1 Mass + TGA"
$ws.Range("A17").Value2 = "1 Mass + F Temp + B Temp"

# ---- Column widths (approximate autofit of the widened columns) ---------
$ws.Columns.Item(1).ColumnWidth = 28.74
$ws.Columns.Item(4).ColumnWidth = 16.17
$ws.Columns.Item(5).ColumnWidth = 23.17
$ws.Columns.Item(6).ColumnWidth = 16.17

$ws.Range("D8").Select()
